$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '67.091.92'
$ws.Range('E2').Value = '  -0.44%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.477.02'
$ws.Range('E3').Value = '  -1.37%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '584.13'
$ws.Range('E5').Value = '  -0.72%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '169.91'
$ws.Range('E6').Value = '  +0.26%  '
$ws.Range('E7').Value = '  +0.11%  '
$ws.Range('E8').Value = '  -1.38%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.476.01'
$ws.Range('E9').Value = '  -1.32%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.137'
$ws.Range('E10').Value = '  +1.35%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.164'
$ws.Range('E11').Value = '  -1.17%  '
$ws.Range('E13').Value = '  -2.90%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '2.925.75'
$ws.Range('E14').Value = '  -2.53%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '25.28'
$ws.Range('E15').Value = '  -2.79%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '66.977.83'
$ws.Range('E16').Value = '  -0.49%  '
$ws.Range('E17').Value = '  -2.08%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.475.85'
$ws.Range('E18').Value = '  -2.99%  '
$ws.Range('E19').Value = '  -6.03%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '350.62'
$ws.Range('E21').Value = '  -3.71%  '
$ws.Range('E22').Value = '  -1.67%  '
$ws.Range('E23').Value = '  +0.57%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '68.56'
$ws.Range('E24').Value = '  -4.38%  '
$ws.Range('E25').Value = '  -6.02%  '
$ws.Range('E26').Value = '  -1.93%  '
$ws.Range('E27').Value = '  -4.83%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.00'
$ws.Range('E28').Value = '  -35.39%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.600.33'
$ws.Range('E29').Value = '  -1.90%  '
$ws.Range('E30').Value = '  -3.96%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '508.27'
$ws.Range('E31').Value = '  -3.69%  '
$ws.Range('E32').Value = '  -6.69%  '
$ws.Range('E33').Value = '  -4.47%  '
$ws.Range('E34').Value = '  -3.44%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.999'
$ws.Range('E35').Value = '  +0.02%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '158.53'
$ws.Range('E36').Value = '  +0.36%  '
$ws.Range('E37').Value = '  -8.98%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '18.23'
$ws.Range('E40').Value = '  -6.48%  '
$ws.Range('E41').Value = '  -0.18%  '
$ws.Range('E42').Value = '  -3.78%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '4.80'
$ws.Range('E43').Value = '  -4.01%  '
$ws.Range('E44').Value = '  -4.27%  '
$ws.Range('E45').Value = '  -2.48%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '38.82'
$ws.Range('E46').Value = '  -1.32%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '142.01'
$ws.Range('E47').Value = '  -3.22%  '
$ws.Range('E48').Value = '  -6.46%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.513'
$ws.Range('E49').Value = '  -5.32%  '
$ws.Range('E50').Value = '  -6.46%  '
$ws.Range('E51').Value = '  -1.40%  '
